# Update "想去人数" (column F) figures on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Mapping of row -> new value for column F on each affected sheet.
$updates = @{
    4  = 865
    6  = 321
    7  = 10040
    8  = 82
    9  = 81
    10 = 136
    17 = 286
    18 = 801
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
